$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2023-11-03 Friday" "2023-11-04 Saturday"

Replace-Text "54×38=2052" "93×24=2232"
Replace-Text "88×40=3520" "33×19=627"
Replace-Text "55×69=3795" "82×89=7298"
Replace-Text "21×45=945" "56×65=3640"
Replace-Text "15×85=1275" "62×30=1860"

Replace-Text "56×44=2464" "38×73=2774"
Replace-Text "38×90=3420" "52×11=572"
Replace-Text "97×72=6984" "80×27=2160"
Replace-Text "74×35=2590" "18×73=1314"
Replace-Text "20×65=1300" "64×92=5888"

Replace-Text "89×37=3293" "67×66=4422"
Replace-Text "86×91=7826" "94×53=4982"
Replace-Text "53×44=2332" "16×90=1440"
Replace-Text "64×13=832" "96×59=5664"
Replace-Text "64×12=768" "98×34=3332"

Replace-Text "29×88=2552" "71×76=5396"
Replace-Text "87×38=3306" "74×51=3774"
Replace-Text "34×20=680" "76×69=5244"
Replace-Text "35×64=2240" "33×56=1848"
Replace-Text "41×13=533" "55×18=990"

Replace-Text "57×11=627" "74×40=2960"
Replace-Text "11×98=1078" "27×89=2403"
Replace-Text "66×25=1650" "58×73=4234"
Replace-Text "62×42=2604" "52×44=2288"
Replace-Text "59×43=2537" "29×38=1102"
